# Applies the edit described by the commit:
#  - In the "Attended by:" paragraph, strike through each attendee's name
#    (Angel, Georgi, Ilia, Mikaeil) while keeping the separating punctuation
#    (", " and the lone ",") unformatted.
#  - Split the former "Planned activities:" paragraph (which also carried the
#    trailing _GoBack bookmark) so that "Planned activities:" and
#    "Highlights:" are each their own paragraph, and the _GoBack bookmark now
#    sits at the end of the "Highlights:" paragraph instead of the
#    "Planned activities:" paragraph.

$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1: "Attended by: Angel, Georgi, Ilia, Mikaeil" ---
# Rebuild the whole paragraph's run sequence in one go (InsertXML on a
# sub-range of a paragraph that doesn't reach the paragraph end can scramble
# run order in this engine, so we always target the *entire* paragraph
# range when reconstructing runs).
$p1 = $d.Paragraphs(1)
$xml1 = "<w:p $w>" +
        "<w:r><w:t>Attended by:</w:t></w:r>" +
        "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
        "<w:r><w:rPr><w:strike/></w:rPr><w:t>Angel</w:t></w:r>" +
        "<w:r><w:t xml:space=`"preserve`">, </w:t></w:r>" +
        "<w:r><w:rPr><w:strike/></w:rPr><w:t>Georgi</w:t></w:r>" +
        "<w:r><w:t xml:space=`"preserve`">, </w:t></w:r>" +
        "<w:r><w:rPr><w:strike/></w:rPr><w:t>Ilia</w:t></w:r>" +
        "<w:r><w:t>,</w:t></w:r>" +
        "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
        "<w:r><w:rPr><w:strike/></w:rPr><w:t>Mikaeil</w:t></w:r>" +
        "</w:p>"
$p1.Range.InsertXML($xml1)

# --- Paragraphs 2 & 3: "Planned activities:" / "Highlights:" ---
# Re-split them and relocate the _GoBack bookmark to the end of the new
# "Highlights:" paragraph.
$p2 = $d.Paragraphs(2)
$p3 = $d.Paragraphs(3)
$combinedRange = $d.Range($p2.Range.Start, $p3.Range.End)
$xml23 = "<w:p $w><w:r><w:t>Planned activities:</w:t></w:r></w:p>" +
         "<w:p $w><w:r><w:t>Highlights:</w:t></w:r>" +
         "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
         "<w:bookmarkEnd w:id=`"0`"/></w:p>"
$combinedRange.InsertXML($xml23)
